$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText($addr, $text) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $text
    $c.ClearFormats()
}

Set-CellText 'D2' '64.405.59'
Set-CellText 'E2' '  -2.55%  '
Set-CellText 'D3' '3.180.58'
Set-CellText 'E3' '  -4.08%  '
Set-CellText 'E4' '  -0.02%  '
Set-CellText 'D5' '571.59'
Set-CellText 'E5' '  -2.36%  '
Set-CellText 'D6' '168.93'
Set-CellText 'E6' '  -7.56%  '
Set-CellText 'D7' '0.609'
Set-CellText 'E7' '  -6.03%  '
Set-CellText 'E8' '  -0.16%  '
Set-CellText 'D9' '3.189.51'
Set-CellText 'E9' '  -3.78%  '
Set-CellText 'E10' '  -3.70%  '
Set-CellText 'D11' '6.79'
Set-CellText 'E11' '  -0.39%  '
Set-CellText 'D12' '0.388'
Set-CellText 'E12' '  -3.12%  '
Set-CellText 'D13' '3.731.51'
Set-CellText 'E13' '  -4.14%  '
Set-CellText 'E14' '  -1.84%  '
Set-CellText 'D15' '64.466.33'
Set-CellText 'E15' '  -2.55%  '
Set-CellText 'E16' '  -3.01%  '
Set-CellText 'E17' '  -3.94%  '
Set-CellText 'D18' '3.189.39'
Set-CellText 'E18' '  -4.57%  '
Set-CellText 'D19' '418.96'
Set-CellText 'E19' '  -1.27%  '
Set-CellText 'D20' '12.98'
Set-CellText 'E20' '  -1.28%  '
Set-CellText 'D21' '5.37'
Set-CellText 'E21' '  -3.08%  '
Set-CellText 'D22' '7.13'
Set-CellText 'E22' '  -3.42%  '
Set-CellText 'D23' '0.999'
Set-CellText 'E23' '  -0.10%  '
Set-CellText 'D24' '70.39'
Set-CellText 'E24' '  -1.91%  '
Set-CellText 'E25' '  -0.16%  '
Set-CellText 'D26' '0.206'
Set-CellText 'E26' '  +2.81%  '
Set-CellText 'D27' '0.490'
Set-CellText 'E27' '  -4.35%  '
Set-CellText 'E28' '  -6.21%  '
Set-CellText 'D29' '8.75'
Set-CellText 'E29' '  -1.63%  '
Set-CellText 'E30' '  -1.15%  '
Set-CellText 'E31' '  -3.34%  '
Set-CellText 'D32' '21.76'
Set-CellText 'E32' '  -2.65%  '
Set-CellText 'D33' '0.999'
Set-CellText 'E33' '  -0.09%  '
Set-CellText 'E34' '  -2.21%  '
Set-CellText 'D35' '6.36'
Set-CellText 'E35' '  -2.88%  '
Set-CellText 'E36' '  -3.71%  '
Set-CellText 'D37' '156.79'
Set-CellText 'E37' '  -2.53%  '
Set-CellText 'E38' '  -4.73%  '
Set-CellText 'D39' '2.707.88'
Set-CellText 'E39' '  -5.94%  '
Set-CellText 'E40' '  -4.91%  '
Set-CellText 'B41' 'EnergySwap'
Set-CellText 'C41' 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-CellText 'D41' '24.30'
Set-CellText 'E41' '  -7.59%  '
Set-CellText 'B42' 'Filecoin'
Set-CellText 'C42' 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-CellText 'D42' '4.22'
Set-CellText 'E42' '  -2.06%  '
Set-CellText 'D43' '39.21'
Set-CellText 'E43' '  -1.83%  '
Set-CellText 'D44' '0.718'
Set-CellText 'E44' '  -5.72%  '
Set-CellText 'D45' '0.0624'
Set-CellText 'E45' '  -5.56%  '
Set-CellText 'D46' '5.58'
Set-CellText 'E46' '  -5.32%  '
Set-CellText 'D47' '0.0265'
Set-CellText 'E47' '  -2.18%  '
Set-CellText 'D48' '292.96'
Set-CellText 'E48' '  -6.65%  '
Set-CellText 'D49' '21.49'
Set-CellText 'E49' '  -6.93%  '
Set-CellText 'B50' 'dogwifhat'
Set-CellText 'C50' 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
Set-CellText 'D50' '2.01'
Set-CellText 'E50' '  -11.77%  '
Set-CellText 'B51' 'FirstDigitalUSD'
Set-CellText 'C51' 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
Set-CellText 'D51' '0.998'
Set-CellText 'E51' '  -0.25%  '
